$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.183.42'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '3.062.97'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '391.90'
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("D6").Value = '101.58'
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").Value = '0.534'
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("D10").Value = '36.73'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '0.0848'
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").Value = '3.544.43'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Value = '18.28'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '7.67'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").Value = '3.061.21'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '0.991'
$ws.Range("E17").Value = '  +2.04%  '
$ws.Range("D18").Value = '10.61'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").Value = '51.152.96'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = '12.25'
$ws.Range("E21").Value = '  -1.62%  '
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = '0.0₃0954'
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").Value = '69.65'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Value = '264.00'
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").Value = '3.14'
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("D26").Value = '7.89'
$ws.Range("E26").Value = '  -6.32%  '
$ws.Range("D27").Value = '26.76'
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = '7.14'
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("D30").Value = '0.163'
$ws.Range("E30").Value = '  -5.19%  '
$ws.Range("D31").Value = '0.105'
$ws.Range("E31").Value = '  -3.10%  '
$ws.Range("D32").Value = '10.46'
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("D33").Value = '0.0487'
$ws.Range("E33").Value = '  +8.76%  '
$ws.Range("D34").Value = '35.76'
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("D35").Value = '2.05'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").Value = '49.94'
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '3.37'
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").Value = '129.63'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '16.53'
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").Value = '0.115'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").Value = '3.78'
$ws.Range("E44").Value = '  +2.35%  '
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '21.71'
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("D47").Value = '2.50'
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").Value = '2.064.84'
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("E50").Value = '  +2.96%  '
$ws.Range("D51").Value = '0.887'
$ws.Range("E51").Value = '  +11.52%  '
